# Revert the "Holds" relation shape that was already implemented differently.
# Removes the "Diamond 32" shape (text "Holds") and its two connectors
# ("Straight Arrow Connector 88" and "Straight Arrow Connector 89") from slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$s.Shapes.Item("Straight Arrow Connector 89").Delete()
$s.Shapes.Item("Straight Arrow Connector 88").Delete()
$s.Shapes.Item("Diamond 32").Delete()
